$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the VLAN range text in cell C3 from "300:500" to "300-500"
$ws.Range("C3").Value = "300-500"

# Reflect the active selection moving to C3 (as seen in the saved file)
$ws.Range("C3").Select()
